# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The stored serial date 45204 (2023-10-05) becomes 45205 (2023-10-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 135; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -eq 45204) {
        $cell.Value = 45205
    }
}
